# Update the last row (2025Q3 / row 29) of the metricas_recorrencia_trimestral sheet
# with refreshed totals, per commit "atualizei dados bibi e add".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C29").Value = 142
$ws.Range("D29").Value = 23
$ws.Range("E29").Value = 119
$ws.Range("F29").Value = 3.95869191049914
